$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7219
$ws1.Range("F7").Value = 183
$ws1.Range("F18").Value = 38
$ws1.Range("F19").Value = 3758
$ws1.Range("F23").Value = 38
$ws1.Range("F26").Value = 2427
$ws1.Range("F28").Value = 299
$ws1.Range("F38").Value = 1449
$ws1.Range("F39").Value = 151

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7219
$ws4.Range("F8").Value = 183
$ws4.Range("F19").Value = 38
$ws4.Range("F20").Value = 3758
$ws4.Range("F24").Value = 38
$ws4.Range("F27").Value = 2427
$ws4.Range("F29").Value = 299
$ws4.Range("F39").Value = 1449
$ws4.Range("F40").Value = 151
